# Generate Report for handoff
# Replaces the stale "57c10ec3-..." handoff-failed entries with a fresh
# "Ready for handoff" status for a newly regenerated source file
# (38b34fba-...md), adds the sibling file that is ready for handoff too
# (ffff554d2bac-...md), and records the per-locale handoff artifact
# (.xlf) + timestamp for each on the zh-cn / de-de detail sheets. The
# ".localization-config" row simply shifts down one row on every sheet.

$wb = $excel.ActiveWorkbook

$baseRepo = "https://github.com/OpenLocalizationTest/oltest/blob"
$srcCommit = "1090a65e6d8fb6a7fd24fddbdcd31aaa016766e1"
$cfgCommit = "042c71ba5ee5a7f37515dd88f1950df56b989302"

$mdFile1 = "38b34fba-f1a2-4a4f-9d37-4a36343c1443.md"
$mdFile2 = "ffff554d2bac-6cc1-4455-958c-6154c9db171d.md"
$cfgFile = ".localization-config"

$xlfZhCn = "38b34fba-f1a2-4a4f-9d37-4a36343c1443.7a436cc468e1d4867240e6964e32e5f9b5dfbbfd.zh-cn.xlf"
$xlfDeDe = "38b34fba-f1a2-4a4f-9d37-4a36343c1443.7a436cc468e1d4867240e6964e32e5f9b5dfbbfd.de-de.xlf"

$readyStatus = "Ready for handoff"
$notLocalized = "Not to be localized"
$epoch = "0001-01-01 00:00:00"

function Set-RowHyperlink($ws, $cellRef, $address, $displayText) {
    # Hyperlinks.Add() on a cell that already owns a hyperlink just stacks a
    # duplicate <hyperlink> entry on top of the stale one in this host, so
    # any pre-existing hyperlink on the row has to be cleared (sheet-wide)
    # before rebuilding it.
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $displayText) | Out-Null
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Clearing any one cell's Hyperlinks resets the whole sheet collection in
# this host, so do it once up front and rebuild every row's link below.
$ws.Range("A2").Hyperlinks.Delete()

Set-RowHyperlink $ws "A2" "$baseRepo/$srcCommit/e2e/$mdFile1" $mdFile1
$ws.Range("B2").Value = $readyStatus
$ws.Range("C2").Value = $readyStatus

Set-RowHyperlink $ws "A3" "$baseRepo/$srcCommit/e2e/$mdFile2" $mdFile2
$ws.Range("B3").Value = $readyStatus
$ws.Range("C3").Value = $readyStatus

Set-RowHyperlink $ws "A4" "$baseRepo/$cfgCommit/$cfgFile" $cfgFile
$ws.Range("B4").Value = $notLocalized
$ws.Range("C4").Value = $notLocalized

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Hyperlinks.Delete()

Set-RowHyperlink $ws "A2" "$baseRepo/$srcCommit/e2e/$mdFile1" $mdFile1
$ws.Range("B2").Value = $readyStatus
Set-RowHyperlink $ws "C2" "$baseRepo/$srcCommit/e2e/$xlfZhCn" $xlfZhCn
$ws.Range("D2").Value = "2016-01-18 12:47:05"
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = "Include"

Set-RowHyperlink $ws "A3" "$baseRepo/$srcCommit/e2e/$mdFile2" $mdFile2
$ws.Range("B3").Value = $readyStatus
Set-RowHyperlink $ws "C3" "$baseRepo/$srcCommit/e2e/$xlfZhCn" $xlfZhCn
$ws.Range("D3").Value = "2016-01-18 12:47:05"
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = "Include"

Set-RowHyperlink $ws "A4" "$baseRepo/$cfgCommit/$cfgFile" $cfgFile
$ws.Range("B4").Value = $notLocalized
$ws.Range("D4").Value = $epoch
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = "Ignored"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Hyperlinks.Delete()

Set-RowHyperlink $ws "A2" "$baseRepo/$srcCommit/e2e/$mdFile1" $mdFile1
$ws.Range("B2").Value = $readyStatus
Set-RowHyperlink $ws "C2" "$baseRepo/$srcCommit/e2e/$xlfDeDe" $xlfDeDe
$ws.Range("D2").Value = "2016-01-18 12:47:14"
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = "Include"

Set-RowHyperlink $ws "A3" "$baseRepo/$srcCommit/e2e/$mdFile2" $mdFile2
$ws.Range("B3").Value = $readyStatus
Set-RowHyperlink $ws "C3" "$baseRepo/$srcCommit/e2e/$xlfDeDe" $xlfDeDe
$ws.Range("D3").Value = "2016-01-18 12:47:14"
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = "Include"

Set-RowHyperlink $ws "A4" "$baseRepo/$cfgCommit/$cfgFile" $cfgFile
$ws.Range("B4").Value = $notLocalized
$ws.Range("D4").Value = $epoch
$ws.Range("G4").Value = $epoch
$ws.Range("H4").Value = "Ignored"

Write-Output "Report regenerated for handoff"
